$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1042134.25
$ws.Cells.Item(17, 10).Value = 1042134.25
$ws.Cells.Item(17, 12).Value = 3126402.75
$ws.Cells.Item(17, 14).Value = -3126738.75

$ws.Cells.Item(40, 8).Value = 1651.0526
$ws.Cells.Item(40, 9).Value = 1833.6364
$ws.Cells.Item(40, 10).Value = 1400
$ws.Cells.Item(40, 11).Value = 1833.6364
$ws.Cells.Item(40, 12).Value = 1400
$ws.Cells.Item(40, 13).Value = -1658.6364
$ws.Cells.Item(40, 14).Value = -1750

$ws.Cells.Item(55, 8).Value = 5187.727
$ws.Cells.Item(55, 9).Value = 1099
$ws.Cells.Item(55, 10).Value = 8595
$ws.Cells.Item(55, 11).Value = 1099
$ws.Cells.Item(55, 12).Value = 8595
$ws.Cells.Item(55, 13).Value = -885
$ws.Cells.Item(55, 14).Value = -9023

$ws.Cells.Item(127, 8).Value = 1174.25
$ws.Cells.Item(127, 9).Value = 648.75
$ws.Cells.Item(127, 10).Value = 1699.75
$ws.Cells.Item(127, 11).Value = 1946.25
$ws.Cells.Item(127, 12).Value = 5099.25
$ws.Cells.Item(127, 13).Value = 3013.75
$ws.Cells.Item(127, 14).Value = -15019.25

$ws.Cells.Item(129, 8).Value = 1032.6078
$ws.Cells.Item(129, 9).Value = 452.42856
$ws.Cells.Item(129, 10).Value = 1124.909
$ws.Cells.Item(129, 11).Value = 1357.28568
$ws.Cells.Item(129, 12).Value = 3374.727
$ws.Cells.Item(129, 13).Value = 3642.71432
$ws.Cells.Item(129, 14).Value = -13374.727

$ws.Cells.Item(138, 8).Value = 2788.203
$ws.Cells.Item(138, 10).Value = 5726.304
$ws.Cells.Item(138, 12).Value = 17178.912
$ws.Cells.Item(138, 14).Value = -27458.912

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 4416.5
$ws.Cells.Item(132, 10).Value = 4699.8
$ws.Cells.Item(132, 12).Value = 14099.4
$ws.Cells.Item(132, 14).Value = -19159.4

$ws.Cells.Item(141, 8).Value = 90000
$ws.Cells.Item(141, 10).Value = 90000
$ws.Cells.Item(141, 12).Value = 90000
$ws.Cells.Item(141, 14).Value = -100360

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3218.6667
$ws.Cells.Item(134, 9).Value = 2937.3333
$ws.Cells.Item(134, 10).Value = 3500
$ws.Cells.Item(134, 11).Value = 8811.999899999999
$ws.Cells.Item(134, 12).Value = 10500
$ws.Cells.Item(134, 13).Value = -6276.999899999999
$ws.Cells.Item(134, 14).Value = -15570

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1773.3334
$ws.Cells.Item(16, 9).Value = 1907.6923
$ws.Cells.Item(16, 10).Value = 900
$ws.Cells.Item(16, 11).Value = 1907.6923
$ws.Cells.Item(16, 12).Value = 900
$ws.Cells.Item(16, 13).Value = -1620.6923
$ws.Cells.Item(16, 14).Value = -1474

$ws.Cells.Item(113, 8).Value = 1773.3334
$ws.Cells.Item(113, 9).Value = 1907.6923
$ws.Cells.Item(113, 10).Value = 900
$ws.Cells.Item(113, 11).Value = 1907.6923
$ws.Cells.Item(113, 12).Value = 900
$ws.Cells.Item(113, 13).Value = 262.3077000000001
$ws.Cells.Item(113, 14).Value = -5240

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(125, 8).Value = 1538.4615
$ws.Cells.Item(125, 9).Value = 800
$ws.Cells.Item(125, 11).Value = 2400
$ws.Cells.Item(125, 13).Value = 2520

$ws.Cells.Item(129, 8).Value = 1602.1351
$ws.Cells.Item(129, 9).Value = 804.8
$ws.Cells.Item(129, 10).Value = 1897.4445
$ws.Cells.Item(129, 11).Value = 2414.4
$ws.Cells.Item(129, 12).Value = 5692.333500000001
$ws.Cells.Item(129, 13).Value = 2585.6
$ws.Cells.Item(129, 14).Value = -15692.3335

$ws.Cells.Item(130, 8).Value = 1395.8334
$ws.Cells.Item(130, 9).Value = 1110
$ws.Cells.Item(130, 11).Value = 3330
$ws.Cells.Item(130, 13).Value = 1690

$ws.Cells.Item(131, 8).Value = 4427.484
$ws.Cells.Item(131, 9).Value = 569.8333
$ws.Cells.Item(131, 10).Value = 5353.32
$ws.Cells.Item(131, 11).Value = 1709.4999
$ws.Cells.Item(131, 12).Value = 16059.96
$ws.Cells.Item(131, 13).Value = 3330.5001
$ws.Cells.Item(131, 14).Value = -26139.96

$ws.Cells.Item(133, 8).Value = 2386
$ws.Cells.Item(133, 9).Value = 2965
$ws.Cells.Item(133, 11).Value = 8895
$ws.Cells.Item(133, 13).Value = -3835

$ws.Cells.Item(134, 8).Value = 1479.0834
$ws.Cells.Item(134, 9).Value = 1107
$ws.Cells.Item(134, 10).Value = 2000
$ws.Cells.Item(134, 11).Value = 3321
$ws.Cells.Item(134, 12).Value = 6000
$ws.Cells.Item(134, 13).Value = 1749
$ws.Cells.Item(134, 14).Value = -16140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(49, 8).Value = 28000
$ws.Cells.Item(49, 10).Value = 28000
$ws.Cells.Item(49, 12).Value = 28000
$ws.Cells.Item(49, 14).Value = -28368

$ws.Cells.Item(70, 8).Value = 5001.34
$ws.Cells.Item(70, 9).Value = 4368.972
$ws.Cells.Item(70, 10).Value = 6627.4287
$ws.Cells.Item(70, 11).Value = 4368.972
$ws.Cells.Item(70, 12).Value = 6627.4287
$ws.Cells.Item(70, 13).Value = -4098.972
$ws.Cells.Item(70, 14).Value = -7167.4287

$ws.Cells.Item(73, 8).Value = 5001.34
$ws.Cells.Item(73, 9).Value = 4368.972
$ws.Cells.Item(73, 10).Value = 6627.4287
$ws.Cells.Item(73, 11).Value = 4368.972
$ws.Cells.Item(73, 12).Value = 6627.4287
$ws.Cells.Item(73, 13).Value = -3432.972
$ws.Cells.Item(73, 14).Value = -8499.4287

$ws.Cells.Item(97, 8).Value = 1172.5
$ws.Cells.Item(97, 9).Value = 1291.25
$ws.Cells.Item(97, 10).Value = 460
$ws.Cells.Item(97, 11).Value = 1291.25
$ws.Cells.Item(97, 12).Value = 460
$ws.Cells.Item(97, 13).Value = -795.25
$ws.Cells.Item(97, 14).Value = -1452

$ws.Cells.Item(140, 8).Value = 39950
$ws.Cells.Item(140, 10).Value = 39950
$ws.Cells.Item(140, 12).Value = 39950
$ws.Cells.Item(140, 14).Value = -50310

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 949827.75
$ws.Cells.Item(2, 10).Value = 2637225
$ws.Cells.Item(2, 12).Value = 2637225
$ws.Cells.Item(2, 14).Value = -2637449

$ws.Cells.Item(132, 8).Value = 6560.793
$ws.Cells.Item(132, 9).Value = 7077
$ws.Cells.Item(132, 10).Value = 5829.5
$ws.Cells.Item(132, 11).Value = 21231
$ws.Cells.Item(132, 12).Value = 17488.5
$ws.Cells.Item(132, 13).Value = -18701
$ws.Cells.Item(132, 14).Value = -22548.5

$ws.Cells.Item(136, 8).Value = 27781812
$ws.Cells.Item(136, 9).Value = 4399.909
$ws.Cells.Item(136, 10).Value = 333333340
$ws.Cells.Item(136, 11).Value = 13199.727
$ws.Cells.Item(136, 12).Value = 1000000020
$ws.Cells.Item(136, 13).Value = -10649.727
$ws.Cells.Item(136, 14).Value = -1000005120

$ws.Cells.Item(141, 8).Value = 53681.11
$ws.Cells.Item(141, 10).Value = 53681.11
$ws.Cells.Item(141, 12).Value = 53681.11
$ws.Cells.Item(141, 14).Value = -64041.11

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).ClearContents()

$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).ClearContents()

$ws.Cells.Item(107, 8).Value = 471.25
$ws.Cells.Item(107, 9).Value = 450
$ws.Cells.Item(107, 10).Value = 484
$ws.Cells.Item(107, 11).Value = 1350
$ws.Cells.Item(107, 12).Value = 1452
$ws.Cells.Item(107, 13).Value = 570
$ws.Cells.Item(107, 14).Value = -5292

$ws.Cells.Item(132, 8).Value = 2153.742
$ws.Cells.Item(132, 9).Value = 1931.8096
$ws.Cells.Item(132, 10).Value = 2619.8
$ws.Cells.Item(132, 11).Value = 5795.4288
$ws.Cells.Item(132, 12).Value = 7859.400000000001
$ws.Cells.Item(132, 13).Value = -3265.4288
$ws.Cells.Item(132, 14).Value = -12919.4

$ws.Cells.Item(136, 8).Value = 2256.125
$ws.Cells.Item(136, 9).Value = 1077.6666
$ws.Cells.Item(136, 10).Value = 3771.2856
$ws.Cells.Item(136, 11).Value = 3232.9998
$ws.Cells.Item(136, 12).Value = 11313.8568
$ws.Cells.Item(136, 13).Value = -682.9998000000001
$ws.Cells.Item(136, 14).Value = -16413.8568

$ws.Cells.Item(141, 8).Value = 181000
$ws.Cells.Item(141, 10).Value = 181000
$ws.Cells.Item(141, 12).Value = 181000
$ws.Cells.Item(141, 14).Value = -191360
